$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    # If the new value looks like a plain number, force the cell to
    # Text format first so Excel keeps storing it as the original
    # literal string (e.g. "26.50", "0.0000175") instead of silently
    # re-typing it as a Number.
    $looksNumeric = $newValue -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $newValue
}

Set-CellValue "D2" '66.968.44'
Set-CellValue "E2" '  +0.74%  '
Set-CellValue "D3" '2.525.36'
Set-CellValue "E3" '  -2.19%  '
Set-CellValue "D5" '589.03'
Set-CellValue "D6" '172.54'
Set-CellValue "E6" '  +3.76%  '
Set-CellValue "E7" '  +0.04%  '
Set-CellValue "D8" '0.526'
Set-CellValue "E8" '  +0.07%  '
Set-CellValue "D9" '2.524.61'
Set-CellValue "E9" '  -2.20%  '
Set-CellValue "E10" '  +0.56%  '
Set-CellValue "E11" '  +2.11%  '
Set-CellValue "D12" '5.14'
Set-CellValue "E12" '  -0.51%  '
Set-CellValue "E13" '  -3.60%  '
Set-CellValue "D14" '26.50'
Set-CellValue "E14" '  -0.60%  '
Set-CellValue "D15" '2.986.42'
Set-CellValue "E15" '  -2.15%  '
Set-CellValue "D16" '0.0000175'
Set-CellValue "E16" '  -1.18%  '
Set-CellValue "D17" '66.901.22'
Set-CellValue "E17" '  +0.73%  '
Set-CellValue "D18" '2.520.27'
Set-CellValue "E18" '  -1.61%  '
Set-CellValue "E19" '  +4.97%  '
Set-CellValue "D20" '11.31'
Set-CellValue "E20" '  -1.02%  '
Set-CellValue "D21" '354.87'
Set-CellValue "E21" '  +0.73%  '
Set-CellValue "E22" '  -1.18%  '
Set-CellValue "D23" '4.60'
Set-CellValue "E23" '  +0.12%  '
Set-CellValue "E24" '  +5.22%  '
Set-CellValue "E25" '  +0.06%  '
Set-CellValue "D26" '69.61'
Set-CellValue "E26" '  +1.28%  '
Set-CellValue "D27" '9.95'
Set-CellValue "E27" '  -0.58%  '
Set-CellValue "D28" '0.999'
Set-CellValue "E28" '  -0.21%  '
Set-CellValue "E29" '  -2.21%  '
Set-CellValue "D30" '0.0₃0974'
Set-CellValue "E30" '  -1.03%  '
Set-CellValue "D31" '531.18'
Set-CellValue "E31" '  -0.48%  '
Set-CellValue "D32" '8.12'
Set-CellValue "E32" '  +1.39%  '
Set-CellValue "D33" '1.32'
Set-CellValue "E33" '  -0.15%  '
Set-CellValue "D34" '1.84'
Set-CellValue "E34" '  -0.37%  '
Set-CellValue "E35" '  -0.91%  '
Set-CellValue "E36" '  +0.01%  '
Set-CellValue "E37" '  -0.23%  '
Set-CellValue "D38" '156.94'
Set-CellValue "E38" '  +0.24%  '
Set-CellValue "E39" '  -0.77%  '
Set-CellValue "E40" '  +1.08%  '
Set-CellValue "D41" '0.353'
Set-CellValue "E41" '  -1.88%  '
Set-CellValue "D43" '5.12'
Set-CellValue "E43" '  +0.25%  '
Set-CellValue "D45" '2.48'
Set-CellValue "E45" '  +3.03%  '
Set-CellValue "D46" '149.05'
Set-CellValue "E46" '  -0.01%  '
Set-CellValue "E47" '  -1.99%  '
Set-CellValue "E48" '  -2.93%  '
Set-CellValue "E49" '  -0.98%  '
Set-CellValue "E50" '  -1.25%  '
Set-CellValue "D51" '0.0758'
